# The xlsx question importer now also filters by "group" (Program area), which
# surfaced a bad row in the national framework questions sheet: row 40 was a
# stray/duplicate entry — it carried the "Flood Monitoring" group label in
# column A, but its Programs/Program description/Question (columns B-D) did
# not belong to that group at all. Delete that entire row; the rows below it
# shift up to close the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$badRow = $ws.Rows.Item(40)
$badRow.Select() | Out-Null
$badRow.Delete()
